$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "MCH232-1"
$ws.Range("C2").Value = "SHOPSTEWARDS/ TU OFFICIALS MANNUAL"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"

# Row 3
$ws.Range("A3").Value = "MCH232-2"
$ws.Range("C3").Value = "CCAWUSA OPEN SCHOOL- SOCIAL THEORY"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"

# Formatting: Calibri 10pt for the new rows (skip column B, which stays untouched)
$fmtA = $ws.Range("A2:A3")
$fmtA.Font.Name = "Calibri"
$fmtA.Font.Size = 10

$fmtCE = $ws.Range("C2:E3")
$fmtCE.Font.Name = "Calibri"
$fmtCE.Font.Size = 10

$fmtGH = $ws.Range("G2:H3")
$fmtGH.Font.Name = "Calibri"
$fmtGH.Font.Size = 10

$fmtF = $ws.Range("F2:F3")
$fmtF.Font.Name = "Calibri"
$fmtF.Font.Size = 10
$fmtF.WrapText = $false
